$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Películas")
$lo = $ws.ListObjects.Item(1)

# Insert a new row inside the sorted table at row 57 (this is where the new
# movie belongs given its score), shifting every row below it down by one.
$ws.Rows.Item(57).Insert()

# Grow the table / autofilter range so it covers the new row too.
$lo.Resize($ws.Range("B2:I97"))

# Fill in the data for the newly released movie.
$ws.Range("B57").Value = "Spenser Confidential"
$ws.Range("C57").Formula = "=AVERAGE(D57,E57,E57,F57,G57,H57,H57,I57)"
$ws.Range("D57").Value = 8
$ws.Range("E57").Value = 8
$ws.Range("F57").Value = 8
$ws.Range("G57").Value = 7
$ws.Range("H57").Value = 6.2
$ws.Range("I57").Value = 5.2

# The newest entry (B57) is highlighted; the previous "latest" entry (B15,
# "Coach Carter") loses the highlight and goes back to the regular style.
$ws.Range("B57").HorizontalAlignment = -4131
$ws.Range("B15").HorizontalAlignment = -4131

$ws.Range("C102").Select()
